$d = $word.ActiveDocument

# The "categoria" table lists Numero | Nome | Contribuicao | Esforco.
# Find the row for "Margarida Ferreira" so we edit the right row even if
# the table is reshuffled, then update her "Esforco" (last) column, which
# currently reads "0 horas", to "3 horas".
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Margarida Ferreira", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$nameCell = $find.Parent.Cells.Item(1)
$t = $nameCell.Range.Tables.Item(1)
$rowIndex = $nameCell.RowIndex
$lastCol = $t.Columns.Count
$targetCell = $t.Cell($rowIndex, $lastCol)

# The cell currently holds "0 horas" as two runs ("0" and " horas"). Replace
# just the leading digit in place so the " horas" run is left untouched.
$cellRange = $targetCell.Range
$digitRange = $d.Range($cellRange.Start, $cellRange.Start + 1)
$digitRange.Text = "3"

# Word tracks the location of the most recent edit with a hidden "_GoBack"
# bookmark. Move it here, right after the new digit and before " horas", to
# match what Word leaves behind after making this edit. Adding a bookmark
# with the same name removes/relocates any bookmark already using that
# name, so the old "_GoBack" (near the end of the document) disappears.
$targetCell2 = $t.Cell($rowIndex, $lastCol)
$markPos = $targetCell2.Range.Start + 1
$bmRange = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
